$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-6 with the new combined tuple-style text values
$ws.Range("A2").Value = "('Forest', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('Island', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A4").Value = "('Mountain', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A5").Value = "('Plains', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A6").Value = "('Swamp', ['Basic Land — Swamp', '({T}: Add {B}.)'])"

# Remove the now-obsolete rows 7 through 31
$ws.Range("A7:A31").ClearContents()
